# Apply updated crypto price/volume figures (columns D and E) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.148.57"
$ws.Range("E2").Value = "  +3.25%  "
$ws.Range("D3").Value = "3.242.39"
$ws.Range("E3").Value = "  +7.14%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'581.71"
$ws.Range("E5").Value = "  +5.27%  "
$ws.Range("D6").Value = "'151.85"
$ws.Range("E6").Value = "  +8.86%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.231.65"
$ws.Range("E8").Value = "  +7.00%  "
$ws.Range("E9").Value = "  +6.39%  "
$ws.Range("D10").Value = "'7.11"
$ws.Range("E10").Value = "  +12.11%  "
$ws.Range("E11").Value = "  +7.46%  "
$ws.Range("D12").Value = "'0.487"
$ws.Range("E12").Value = "  +6.11%  "
$ws.Range("D13").Value = "'37.79"
$ws.Range("E13").Value = "  +4.70%  "
$ws.Range("D14").Value = "'0.0000234"
$ws.Range("E14").Value = "  +7.29%  "
$ws.Range("D15").Value = "3.767.49"
$ws.Range("E15").Value = "  +7.37%  "
$ws.Range("D16").Value = "66.202.24"
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("D17").Value = "'548.77"
$ws.Range("E17").Value = "  +14.44%  "
$ws.Range("D18").Value = "3.250.51"
$ws.Range("E18").Value = "  +7.37%  "
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("D20").Value = "'7.12"
$ws.Range("E20").Value = "  +7.31%  "
$ws.Range("D21").Value = "'14.53"
$ws.Range("E21").Value = "  +7.52%  "
$ws.Range("D22").Value = "'0.744"
$ws.Range("E22").Value = "  +9.44%  "
$ws.Range("D23").Value = "'7.85"
$ws.Range("E23").Value = "  +11.41%  "
$ws.Range("D24").Value = "'13.47"
$ws.Range("E24").Value = "  +8.01%  "
$ws.Range("D25").Value = "'81.29"
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'9.35"
$ws.Range("E27").Value = "  +20.17%  "
$ws.Range("D28").Value = "'2.98"
$ws.Range("E28").Value = "  +10.19%  "
$ws.Range("E29").Value = "  +7.33%  "
$ws.Range("D30").Value = "'27.75"
$ws.Range("E30").Value = "  +8.08%  "
$ws.Range("D31").Value = "'2.76"
$ws.Range("E31").Value = "  +7.05%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "'1.17"
$ws.Range("E33").Value = "  +6.19%  "
$ws.Range("D34").Value = "'567.02"
$ws.Range("E34").Value = "  +9.66%  "
$ws.Range("D35").Value = "'5.67"
$ws.Range("E35").Value = "  +5.46%  "
$ws.Range("D36").Value = "'6.36"
$ws.Range("E36").Value = "  +7.76%  "
$ws.Range("D37").Value = "'55.20"
$ws.Range("E37").Value = "  +5.62%  "
$ws.Range("D38").Value = "'0.0453"
$ws.Range("E38").Value = "  +13.68%  "
$ws.Range("D39").Value = "'0.0862"
$ws.Range("E39").Value = "  +8.85%  "
$ws.Range("D40").Value = "'0.130"
$ws.Range("E40").Value = "  +7.55%  "
$ws.Range("D41").Value = "'2.97"
$ws.Range("E41").Value = "  +9.99%  "
$ws.Range("D42").Value = "3.201.09"
$ws.Range("E42").Value = "  +11.49%  "
$ws.Range("D43").Value = "'8.62"
$ws.Range("E43").Value = "  +4.73%  "
$ws.Range("D44").Value = "'0.284"
$ws.Range("E44").Value = "  +17.92%  "
$ws.Range("E45").Value = "  +12.28%  "
$ws.Range("D46").Value = "'26.50"
$ws.Range("E46").Value = "  +6.95%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "0.0₃0558"
$ws.Range("E48").Value = "  +6.00%  "
$ws.Range("D49").Value = "'125.76"
$ws.Range("E49").Value = "  +4.88%  "
$ws.Range("E50").Value = "  +4.76%  "
$ws.Range("D51").Value = "'2.21"
$ws.Range("E51").Value = "  +9.99%  "
